# Applies the Leve-profit recalculation updates captured in the commit diff.
# Each row below mirrors one <row> hunk from the OOXML diff: columns H-N hold
# the recomputed price/profit figures for that leve row. Cells absent from the
# "after" XML are cleared with ClearContents() rather than set to 0.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3678.8
$ws.Range("I2").Value = 1750.5
$ws.Range("K2").Value = 1750.5
$ws.Range("M2").Value = -1637.5
$ws.Range("H43").Value = 2748
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2748
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2748
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2886
$ws.Range("H70").Value = 3616.25
$ws.Range("I70").Value = 2333
$ws.Range("J70").Value = 4044
$ws.Range("K70").Value = 6999
$ws.Range("L70").Value = 12132
$ws.Range("M70").Value = -6729
$ws.Range("N70").Value = -12672
$ws.Range("H73").Value = 3616.25
$ws.Range("I73").Value = 2333
$ws.Range("J73").Value = 4044
$ws.Range("K73").Value = 6999
$ws.Range("L73").Value = 12132
$ws.Range("M73").Value = -6063
$ws.Range("N73").Value = -14004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 4677.25
$ws.Range("I31").Value = 4677.25
$ws.Range("K31").Value = 4677.25
$ws.Range("M31").Value = -4383.25
$ws.Range("H45").Value = 1831.25
$ws.Range("I45").Value = 1815.909
$ws.Range("K45").Value = 1815.909
$ws.Range("M45").Value = -1438.909

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 5558.4287
$ws.Range("J88").Value = 5558.4287
$ws.Range("L88").Value = 5558.4287
$ws.Range("N88").Value = -6370.4287
$ws.Range("H91").Value = 5558.4287
$ws.Range("J91").Value = 5558.4287
$ws.Range("L91").Value = 5558.4287
$ws.Range("N91").Value = -8366.4287
$ws.Range("H100").Value = 41749.832
$ws.Range("J100").Value = 41749.832
$ws.Range("L100").Value = 41749.832
$ws.Range("N100").Value = -43913.832
$ws.Range("H134").Value = 5753.533
$ws.Range("I134").Value = 5753.533
$ws.Range("K134").Value = 17260.599
$ws.Range("M134").Value = -14725.599

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2490
$ws.Range("I31").Value = 2490
$ws.Range("K31").Value = 2490
$ws.Range("M31").Value = -2195
$ws.Range("H34").Value = 2490
$ws.Range("I34").Value = 2490
$ws.Range("K34").Value = 2490
$ws.Range("M34").Value = -2288
$ws.Range("H41").Value = 18250
$ws.Range("H62").Value = 5005
$ws.Range("I62").Value = 5005
$ws.Range("K62").Value = 5005
$ws.Range("M62").Value = -4381
$ws.Range("H65").Value = 5005
$ws.Range("I65").Value = 5005
$ws.Range("K65").Value = 25025
$ws.Range("M65").Value = -21905
$ws.Range("H104").Value = 69275
$ws.Range("J104").Value = 69275
$ws.Range("L104").Value = 69275
$ws.Range("N104").Value = -74517
$ws.Range("H132").Value = 3747.5
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 3995
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 11985
$ws.Range("M132").Value = -7970
$ws.Range("N132").Value = -17045

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H134").Value = 982.5
$ws.Range("I134").Value = 982.5
$ws.Range("K134").Value = 2947.5
$ws.Range("M134").Value = 2122.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9618517
$ws.Range("I122").Value = 12502413
$ws.Range("K122").Value = 37507239
$ws.Range("M122").Value = -37504789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3931
$ws.Range("I7").Value = 3901.5
$ws.Range("K7").Value = 3901.5
$ws.Range("M7").Value = -3789.5
$ws.Range("H22").Value = 3391.8
$ws.Range("J22").Value = 3489.75
$ws.Range("L22").Value = 3489.75
$ws.Range("N22").Value = -4079.75
$ws.Range("H27").Value = 3391.8
$ws.Range("J27").Value = 3489.75
$ws.Range("L27").Value = 3489.75
$ws.Range("N27").Value = -3703.75
$ws.Range("H46").Value = 2212.4119
$ws.Range("I46").Value = 2269.5
$ws.Range("K46").Value = 2269.5
$ws.Range("M46").Value = -2081.5
$ws.Range("H82").Value = 1535.25
$ws.Range("I82").Value = 1270
$ws.Range("K82").Value = 1270
$ws.Range("M82").Value = -909
$ws.Range("H85").Value = 1535.25
$ws.Range("I85").Value = 1270
$ws.Range("K85").Value = 1270
$ws.Range("M85").Value = -22
$ws.Range("H122").Value = 3501.5
$ws.Range("I122").Value = 3501.5
$ws.Range("K122").Value = 10504.5
$ws.Range("M122").Value = -8054.5
$ws.Range("H126").Value = 3931
$ws.Range("I126").Value = 3901.5
$ws.Range("K126").Value = 11704.5
$ws.Range("M126").Value = -9234.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 21834.2
$ws.Range("J69").Value = 21834.2
$ws.Range("L69").Value = 21834.2
$ws.Range("N69").Value = -23332.2
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H72").Value = 21834.2
$ws.Range("J72").Value = 21834.2
$ws.Range("L72").Value = 65502.60000000001
$ws.Range("N72").Value = -72990.60000000001
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H81").Value = 800
$ws.Range("I81").Value = 800
$ws.Range("K81").Value = 1600
$ws.Range("M81").Value = -539
$ws.Range("H84").Value = 800
$ws.Range("I84").Value = 800
$ws.Range("K84").Value = 8000
$ws.Range("M84").Value = -2696
$ws.Range("H122").Value = 2362
$ws.Range("I122").Value = 2249
$ws.Range("J122").Value = 2475
$ws.Range("K122").Value = 6747
$ws.Range("L122").Value = 7425
$ws.Range("M122").Value = -4297
$ws.Range("N122").Value = -12325
$ws.Range("H132").Value = 3399.5
$ws.Range("I132").Value = 2599.5
$ws.Range("K132").Value = 7798.5
$ws.Range("M132").Value = -5268.5
$ws.Range("H136").Value = 1386.3846
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

